$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1. Update the document creation date in the metadata table.
Replace-Text "2020-09-03" "2020-09-07"

# 2. Update the work-programme bullet list (two bullets).
Replace-Text `
    "establishing a Smart Energy Observatory of ~10,000 smart metered households across Great Britain (GB);" `
    "establishing an Observatory panel of ~10,000 smart metered households across Great Britain (GB);"

Replace-Text `
    "developing a mechanism for other researchers to access smart meter data for other household samples from whom they have obtained informed consent." `
    "provisioning of Observatory data to UK researchers via a secure lab environment."

# 3. Split the intro paragraph: trim the first paragraph and push the
#    "Smart meter data collection..." sentence into its own new paragraph.
Replace-Text `
    "(Data Access section for more details). Smart meter data collection start dates vary by participant; the earliest start date is in August 2018. This data release contains all data available up until 31st July 2020. Future data releases will include more recent data and data from participants recruited in the second and third recruitment waves." `
    "(Data Access section for more details).^pSmart meter data collection start dates vary by participant; the earliest start date is in August 2018. This data release contains all data available up until 31st July 2020. Future data releases will include more recent data and data from participants recruited in the second and third waves of participant recruitment."

# The new paragraph inherits the "FirstParagraph" style from the split point;
# fix it up to be a regular "Body Text" paragraph like the rest of the section.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Smart meter data collection start dates vary by participant*") {
        $p.Style = "Body Text"
    }
}

# 4. Tweak the wording of the first run in the PUPRN paragraph and the run
#    describing the identifier.
Replace-Text `
    "The datasets can be linked at the household level using the" `
    "The datasets can be linked at the household level using the pseudo-anonymous"

Replace-Text `
    "pseudo-anonymous identifier included in each dataset. The climate data can be linked through the" `
    "(pseudo-UPRN (Unique Property Reference Number)) identifier included in each dataset. The climate data can be linked through the"

# 5. Smart meter data section wording tweaks.
Replace-Text `
    "In order to participate a household must have a DCC-enrolled electricity smart meter (SMETS2 or upgraded SMETS1)." `
    "In order to participate a household must have a DCC-enrolled electricity smart meter (SMETS2 or DCC-enrolled SMETS1)."

Replace-Text `
    "Currently 1612 participants have daily smart meter reads in the datset and 1688 have half-hourly reads in the dataset. 427 participants do not have a gas meter (that we are able to access)." `
    "As of 2020-07-31 1612 participants have daily smart meter reads in the datset and 1688 have half-hourly reads in the dataset. 427 participants do not have a gas meter (that we are able to access)."

# 6. Fill in actual row/column counts that had been unrendered inline R code.
Replace-Text "r ncol_d fields" "22 fields"
Replace-Text "r nrow_hh records" "37049124 records"
Replace-Text "r ncol_hh fields" "20 fields"
Replace-Text "r nrow_reads records" "7897 records"
Replace-Text "r ncol_reads fields" "25 fields"
Replace-Text "r nrow_pp records" "1708 records"
Replace-Text "r ncol_pp fields" "39 fields"

# 7. Survey completion wording + field count.
Replace-Text `
    "Survey data exists for 1673 participants, of whom 1313 completed the survey." `
    "Survey data exists for 1673 participants, of whom 1313 completed the survey in full."

Replace-Text "151 fields" "155 fields"
